$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 24, shifting existing rows 24+ down to 26+
$ws.Rows.Item(24).Resize(2).Insert()

# Fill the newly inserted row 24 (copy of old row 24 "Especial" entry with updated date/price)
$ws.Cells.Item(24, 1).Value = 1
$ws.Cells.Item(24, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(24, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(24, 4).Value = 44459
$ws.Cells.Item(24, 5).Value = 15
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100108
$ws.Cells.Item(24, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(24, 9).Value = 100108002
$ws.Cells.Item(24, 10).Value = "Mango"
$ws.Cells.Item(24, 11).Value = "Sin especificar"
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 450
$ws.Cells.Item(24, 14).Value = 7000
$ws.Cells.Item(24, 15).Value = 8000
$ws.Cells.Item(24, 16).Value = 7500
$ws.Cells.Item(24, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(24, 18).Value = "Perú"
$ws.Cells.Item(24, 19).Value = 1875
$ws.Cells.Item(24, 20).Value = 4

# Fill the newly inserted row 25 (copy of old row 25 "Primera" entry with updated date/price)
$ws.Cells.Item(25, 1).Value = 1
$ws.Cells.Item(25, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(25, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(25, 4).Value = 44459
$ws.Cells.Item(25, 5).Value = 15
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100108
$ws.Cells.Item(25, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(25, 9).Value = 100108002
$ws.Cells.Item(25, 10).Value = "Mango"
$ws.Cells.Item(25, 11).Value = "Sin especificar"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 600
$ws.Cells.Item(25, 14).Value = 7000
$ws.Cells.Item(25, 15).Value = 8000
$ws.Cells.Item(25, 16).Value = 7500
$ws.Cells.Item(25, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(25, 18).Value = "Perú"
$ws.Cells.Item(25, 19).Value = 1875
$ws.Cells.Item(25, 20).Value = 4
